# updated results with new hyperparameter
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mnist")

# ------------------------------------------------------------------
# 1) FPR columns (I = baard_2stage, J = baard_3stage) for every data
#    row in the table now report the same re-run value (3.8) instead
#    of the old per-row 6.3 / 8.6 pair. Rows 4-22 (row 13 excluded
#    below, handled separately) all get this update.
# ------------------------------------------------------------------
$fprRows = @(4,5,6,7,8,9,10,11,12,14,15,16,17,18,19,20,21,22)
foreach ($r in $fprRows) {
    $ws.Cells.Item($r, 9).Value = 3.8   # column I
    $ws.Cells.Item($r, 10).Value = 3.8  # column J
}

# ------------------------------------------------------------------
# 2) Row 13 ("boundary", Adv_param column B) lost its adversarial
#    parameter value - the cell is cleared back to blank.
# ------------------------------------------------------------------
$ws.Range("B13").ClearContents()

# ------------------------------------------------------------------
# 3) Row 23 ("line" attack, first sub-row) now reports fresh numbers
#    from the rerun with the new hyperparameter.
# ------------------------------------------------------------------
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 99.90000000000001
$ws.Range("D23").Value = 99.90000000000001
$ws.Range("E23").Value = 99.09999999999999
$ws.Range("F23").Value = 96.8
$ws.Range("G23").Value = 98.40000000000001
$ws.Range("H23").Value = 92.5
$ws.Range("I23").Value = 3.8
$ws.Range("J23").Value = 3.8
$ws.Range("L23").Value = 0.1

# ------------------------------------------------------------------
# 4) Row 24 becomes the start of the "watermark" attack block (it used
#    to start at row 26) - set the label and the refreshed figures.
# ------------------------------------------------------------------
$ws.Range("A24").Value = "watermark"
$ws.Range("B24").Value = 0.3
$ws.Range("C24").Value = 98.8
$ws.Range("D24").Value = 98.8
$ws.Range("F24").Value = 98.40000000000001
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 97.40000000000001
$ws.Range("I24").Value = 3.8
$ws.Range("J24").Value = 3.8

# ------------------------------------------------------------------
# 5) Row 25 (second/last "watermark" sub-row) also gets refreshed
#    figures (it used to hold what is now row 23's old data).
# ------------------------------------------------------------------
$ws.Range("B25").Value = 0.6
$ws.Range("C25").Value = 95
$ws.Range("D25").Value = 95
$ws.Range("E25").Value = 98.40000000000001
$ws.Range("F25").Value = 93.7
$ws.Range("G25").Value = 93.59999999999999
$ws.Range("H25").Value = 91.7
$ws.Range("I25").Value = 3.8
$ws.Range("J25").Value = 3.8
$ws.Range("L25").Value = 0

# ------------------------------------------------------------------
# 6) The old rows 26-27 (second half of "watermark") are gone - the
#    table now ends at row 25. Re-merge A24:A25 for the "watermark"
#    label (it previously spanned A26:A27), and drop the old A23:A25
#    merge since "line" (row 23) is no longer part of a merged block.
# ------------------------------------------------------------------
$ws.Range("A23:A25").UnMerge()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(26).Delete()
$ws.Range("A24:A25").Merge()
